# ------------------------------------------------------------------
# Edit: "Added documentation of the two design reviews."
#  1. Tidy the "Cornell Cup Final Report" bullet: merge the split
#     "Report" run (removing the gramStart/gramEnd proof markers) and
#     drop the trailing period.
#  2. Drop trailing periods from three more bullets (Graphical
#     Timeline Appendix, Timeline Appendix, Bill Of Materials
#     Appendix).
#  3. Insert two new bullets documenting the First and Second Design
#     Reviews right after the "Bill Of Materials Appendix" bullet.
#  4. Move the _GoBack bookmark from the end of the document to the
#     end of the (now-first) "Cornell Cup Final Report" bullet.
# ------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Cornell Cup Final Report bullet -----------------------------------
# Collapses the "Cornell Cup Final "/"Report"/" :"/" Word document..." runs
# (and the gramStart/gramEnd proofErr markers between them) into one run,
# and removes the trailing period at the same time.
$d.Content.Find.Execute(
    "Cornell Cup Final Report : Word document containing final report.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Cornell Cup Final Report : Word document containing final report", 2) | Out-Null

# Re-split so only "Cornell Cup Final Report" stays bold.
$p1 = $d.Paragraphs.Item(3)
$boldLen = ("Cornell Cup Final Report").Length
$nonBoldStart = $p1.Range.Start + $boldLen
$nonBoldEnd = $p1.Range.End - 1
$nonBoldRange = $d.Range($nonBoldStart, $nonBoldEnd)
$nonBoldRange.Font.Bold = 0

# --- 2. Drop trailing periods on three more bullets ------------------------
$d.Content.Find.Execute(
    "visual representation of the timeline.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "visual representation of the timeline", 2) | Out-Null

$d.Content.Find.Execute(
    "listing of dates for the timeline.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "listing of dates for the timeline", 2) | Out-Null

$d.Content.Find.Execute(
    "bought and used for the project.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "bought and used for the project", 2) | Out-Null

# --- 3. Insert the two new Design Review bullets ----------------------------
# "Bill Of Materials Appendix" is paragraph 6; insert two cloned
# list-paragraphs right after it (InsertParagraphAfter duplicates the
# bold ListParagraph formatting of paragraph 6). For each new bullet the
# heading ("First/Second Design Review: ") is explicitly re-bolded and the
# trailing description explicitly un-bolded, rather than relying on
# inherited run formatting (which only survives for the very first clone).
$bom = $d.Paragraphs.Item(6)
$bom.Range.InsertParagraphAfter() | Out-Null

$firstReview = $d.Paragraphs.Item(7)
$firstReview.Range.InsertAfter("First Design Review: PDF file documenting how our presentation was evaluated and comments suggested") | Out-Null
$boldLen = ("First Design Review: ").Length
$boldStart = $firstReview.Range.Start
$boldEnd = $firstReview.Range.Start + $boldLen
$d.Range($boldStart, $boldEnd).Font.Bold = 1
$nonBoldStart = $boldEnd
$nonBoldEnd = $firstReview.Range.End - 1
$d.Range($nonBoldStart, $nonBoldEnd).Font.Bold = 0

$firstReview.Range.InsertParagraphAfter() | Out-Null
$secondReview = $d.Paragraphs.Item(8)
$secondReview.Range.InsertAfter("Second Design Review: Word document explaining main talking points of the design review") | Out-Null
$boldLen2 = ("Second Design Review: ").Length
$boldStart2 = $secondReview.Range.Start
$boldEnd2 = $secondReview.Range.Start + $boldLen2
$d.Range($boldStart2, $boldEnd2).Font.Bold = 1
$nonBoldStart2 = $boldEnd2
$nonBoldEnd2 = $secondReview.Range.End - 1
$d.Range($nonBoldStart2, $nonBoldEnd2).Font.Bold = 0

# --- 4. Move the _GoBack bookmark to the end of the Cornell Cup bullet -----
# (Bookmarks.Add with a range collapsed exactly at a paragraph's final
# character position lands incorrectly, so pad with a throw-away
# character, anchor the bookmark before it, then delete the character.)
$cornell = $d.Paragraphs.Item(3)
$endPos = $cornell.Range.End - 1
$pad = $d.Range($endPos, $endPos)
$pad.InsertAfter("X") | Out-Null
$bookmarkRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
$d.Range($endPos, $endPos + 1).Delete() | Out-Null

Write-Output "done"
